$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.068.86"
$ws.Range("E2").Value = "  +1.83%  "

$ws.Range("D3").Value = "3.830.93"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'635.82"
$ws.Range("E5").Value = "  +5.60%  "

$ws.Range("D6").Value = "'165.42"
$ws.Range("E6").Value = "  -0.56%  "

$ws.Range("D7").Value = "3.828.82"
$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("D11").Value = "'0.454"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").Value = "'6.69"
$ws.Range("E12").Value = "  +4.91%  "

$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("D14").Value = "'35.97"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").Value = "4.473.90"
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("D16").Value = "3.839.04"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("D17").Value = "69.112.67"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").Value = "'18.26"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "'466.85"
$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("D22").Value = "'9.76"
$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("E23").Value = "  +0.75%  "

$ws.Range("E24").Value = "  +2.67%  "

$ws.Range("D25").Value = "'83.92"
$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").Value = "'12.06"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "3.981.09"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").Value = "'2.69"
$ws.Range("E31").Value = "  -3.34%  "

$ws.Range("D32").Value = "'2.26"
$ws.Range("E32").Value = "  +1.82%  "

$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("D34").Value = "'29.20"
$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("D35").Value = "'9.14"
$ws.Range("E35").Value = "  +0.62%  "

$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("E37").Value = "  +1.91%  "

$ws.Range("E38").Value = "  +7.03%  "

$ws.Range("D39").Value = "'5.94"
$ws.Range("E39").Value = "  +2.15%  "

$ws.Range("D40").Value = "'3.29"
$ws.Range("E40").Value = "  +1.13%  "

$ws.Range("E41").Value = "  -1.73%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").Value = "'156.82"
$ws.Range("E44").Value = "  +3.62%  "

$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("D46").Value = "'1.43"
$ws.Range("E46").Value = "  +3.37%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'46.86"
$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'42.91"
$ws.Range("E48").Value = "  -5.33%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.43"
$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.89"
$ws.Range("E50").Value = "  +2.30%  "

$ws.Range("D51").Value = "'382.49"
$ws.Range("E51").Value = "  -2.47%  "
